# Power Models and Textures.xlsx
# Insert a new "textures\water_jungle.igb" row into Table2 on Sheet1, in its
# correct alphabetically-sorted position (between "water_insect" and
# "water_monument", i.e. before the existing row 139), shifting every row
# below it down by one and growing the table/formatting ranges accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Make room: insert a blank row at 139, shifting 139-148 -> 140-149 ---
$ws.Rows.Item(139).Insert()

# --- 2. Grow the table ("Table2") to cover the new row ---
$tbl = $ws.ListObjects.Item("Table2")
$tbl.Resize($ws.Range("A1:G149"))

# --- 3. Fill in the new row's data (same PC/GC/PS2/PSP/Xbox source + Notes
#        pattern as its alphabetical neighbours, e.g. "water_monument") ---
$ws.Range("A139").Value = "textures\water_jungle.igb"
$ws.Range("B139").Value = "XML2 PC"
$ws.Range("C139").Value = "XML2 GameCube"
$ws.Range("D139").Value = "XML2 PS2"
$ws.Range("E139").Value = "XML2 PSP"
$ws.Range("F139").Value = "XML2 Xbox"
$ws.Range("G139").Value = "1b. Used in all versions of XML2 but not originally in permanent"

# --- 4. Extend every conditional-formatting rule's range by one row
#        (A2:A148 -> A2:A149, B2:B148 -> B2:B149, B2:F148 -> B2:F149, ...)
#        Read every rule's current AppliesTo column span first (a single
#        query over the whole table returns every distinct rule exactly
#        once), then apply the resize, so rules aren't revisited/doubled. ---
$allFC = $ws.Range("A2:G149").FormatConditions
$count = $allFC.Count
$spans = @()
for ($i = 1; $i -le $count; $i++) {
    $applies = $allFC.Item($i).AppliesTo
    $spans += ,@($applies.Column, $applies.Columns.Count)
}
for ($i = 1; $i -le $count; $i++) {
    $col1 = $spans[$i - 1][0]
    $colN = $spans[$i - 1][1]
    $newRange = $ws.Range($ws.Cells.Item(2, $col1), $ws.Cells.Item(149, $col1 + $colN - 1))
    $ws.Range("A2:G149").FormatConditions.Item($i).ModifyAppliesToRange($newRange)
}

# --- 5. Move the on-screen selection to match (A149 -> A135) ---
$null = $ws.Range("A135").Select()
